$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix A96: the date/time value was corrected ---
$ws.Range("A96").Value = 45481.2916666667

# --- Add new row 97 with the latest data point ---
$ws.Range("A97").Value = 45482.5247569444

# Match the date/time style already used by column A (same numFmt + font)
# instead of letting Excel invent a brand new cell style.
$ws.Range("A96").Copy()
$ws.Range("A97").PasteSpecial(-4122)
$excel.CutCopyMode = 0

$ws.Range("B97").Value = 13500
$ws.Range("C97").Value = 3.50999999046326
$ws.Range("D97").Value = 3.13000011444092
$ws.Range("E97").Value = 3.44000005722046
$ws.Range("F97").Value = 3.1800000667572

# G97 (adj_close) holds the close price formatted as text, matching the
# pattern used by every other row in this column. A direct string
# assignment would be auto-coerced back into a number by Excel, so we
# build it as a text formula result first and then freeze it to a plain
# value via copy / paste-special, which keeps it tagged as text without
# touching the cell's number format/style.
$ws.Range("G97").Formula = '=""&"3.1800000667572"'
$ws.Range("G97").Copy()
$ws.Range("G97").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("H97").Value = "ESPE.MI"
